$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: row -> (Nome completo, CPF)
$data = @{
    2  = @("Isaac Moreira", "149.578.326-09")
    3  = @("Igor Barros", "709.815.623-02")
    4  = @("Vitória Almeida", "160.478.932-87")
    5  = @("Amanda Fernandes", "134.568.739-25")
    6  = @("Lucas Pereira", "268.235.174-60")
    7  = @("Samuel Costa", "874.123.756-50")
    8  = @("Pedro Almeida", "194.854.362-91")
    9  = @("João Vitor", "584.963.217-49")
    10 = @("Beatriz Gomes", "953.142.689-01")
    11 = @("Leonardo Lima", "362.497.528-09")
    12 = @("Gustavo Batista", "851.324.178-12")
    13 = @("Daniel Cunha", "342.178.543-65")
    14 = @("Matheus Barros", "146.983.254-17")
    15 = @("Rafael Oliveira", "734.862.319-05")
    16 = @("Renan Cardoso", "694.251.784-23")
    17 = @("Diego Vieira", "194.825.763-02")
    18 = @("Bruno Mendes", "563.712.358-46")
    19 = @("Thiago Barbosa", "167.925.138-80")
    20 = @("Henrique Ferreira", "643.187.925-48")
    21 = @("Vinícius Oliveira", "694.871.352-10")
    22 = @("André Dias", "712.358.194-83")
    23 = @("Nicole Borges", "157.923.485-06")
    24 = @("Lucas Cavalcante", "189.725.143-40")
    25 = @("Sérgio Correia", "583.491.763-90")
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
}
